$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily data row (10/16/2025) at the bottom of the table,
# mirroring the plain-text date formatting and default (no) cell style
# used by the existing data rows.
$row = 45

$ws.Cells.Item($row, 1).Value = "'10/16/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.1765432923704638
$ws.Cells.Item($row, 3).Value = 0.8234567076295362
